$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31; existing rows 31-36 shift down to 32-37.
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new record.
$ws.Range("A31").Value = 5
$ws.Range("B31").Value = "Macroferia Regional de Talca"
$ws.Range("C31").Value = "Maule"
$ws.Range("D31").Value = 44769
$ws.Range("E31").Value = 7
$ws.Range("F31").Value = 100112040
$ws.Range("G31").Value = "Cilantro"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 200
$ws.Range("K31").Value = 12000
$ws.Range("L31").Value = 12000
$ws.Range("M31").Value = 12000
$ws.Range("N31").Value = "$/caja 36 atados"
$ws.Range("O31").Value = "Región del Maule"
$ws.Range("P31").Value = 333
$ws.Range("Q31").Value = 36
$ws.Range("R31").Value = "Hortaliza"
